$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain-text values such as
# "309.16" or "1.93%". Mark these columns as Text format before writing the
# new values so Excel keeps them as literal strings instead of silently
# re-interpreting them as numbers / percentages.
$dRange = $ws.Range("D2:D51")
$eRange = $ws.Range("E2:E51")
$dRange.NumberFormat = "@"
$eRange.NumberFormat = "@"

# Updated price / volume figures for this data refresh
$ws.Range("D2").Value = "309.16"
$ws.Range("E2").Value = "1.93%"
$ws.Range("D3").Value = "38.91"
$ws.Range("E3").Value = "9.06%"
$ws.Range("D4").Value = "5.091"
$ws.Range("E4").Value = "1.18%"
$ws.Range("D5").Value = "0.08193"
$ws.Range("E5").Value = "3.72%"
$ws.Range("D6").Value = "2.033"
$ws.Range("E6").Value = "10.08%"
$ws.Range("D7").Value = "7.907"
$ws.Range("E7").Value = "1.58%"
$ws.Range("E8").Value = "1.02%"
$ws.Range("D9").Value = "0.1423"
$ws.Range("E9").Value = "5.40%"
$ws.Range("D10").Value = "0.1950"
$ws.Range("E10").Value = "3.06%"
$ws.Range("D11").Value = "0.09246"
$ws.Range("E11").Value = "2.13%"
$ws.Range("E12").Value = "-0.13%"
$ws.Range("E13").Value = "0.26%"
$ws.Range("D14").Value = "0.001405"
$ws.Range("E14").Value = "0.21%"
$ws.Range("D15").Value = "0.005887"
$ws.Range("E15").Value = "-3.95%"
$ws.Range("D16").Value = "3.788"
$ws.Range("E16").Value = "1.80%"
$ws.Range("D17").Value = "4.176"
$ws.Range("E17").Value = "1.88%"
$ws.Range("D18").Value = "3.441"
$ws.Range("E18").Value = "4.33%"
$ws.Range("D19").Value = "0.3452"
$ws.Range("E19").Value = "0.40%"
$ws.Range("D20").Value = "0.1302"
$ws.Range("E20").Value = "-2.99%"
$ws.Range("D21").Value = "4.837"
$ws.Range("E21").Value = "-6.42%"
$ws.Range("D22").Value = "0.2352"
$ws.Range("E22").Value = "7.30%"
$ws.Range("D23").Value = "0.04466"
$ws.Range("E23").Value = "1.31%"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").Value = "0.46%"
$ws.Range("E25").Value = "-9.62%"
$ws.Range("D27").Value = "0.0001301"
$ws.Range("E27").Value = "0.04%"
$ws.Range("D39").Value = "0.02123"
$ws.Range("E39").Value = "9.51%"
$ws.Range("D40").Value = "0.05175"
$ws.Range("E40").Value = "0.95%"
$ws.Range("D41").Value = "0.007470"
$ws.Range("E41").Value = "-1.96%"
$ws.Range("D42").Value = "0.01012"
$ws.Range("E42").Value = "-0.29%"
$ws.Range("D43").Value = "0.1368"
$ws.Range("E43").Value = "1.92%"
$ws.Range("E44").Value = "-1.81%"
$ws.Range("D45").Value = "0.009677"
$ws.Range("E45").Value = "-5.13%"
$ws.Range("D46").Value = "0.00006308"
$ws.Range("E46").Value = "2.68%"
$ws.Range("E47").Value = "0.03%"
$ws.Range("E48").Value = "-0.24%"
$ws.Range("E49").Value = "-3.51%"
$ws.Range("E50").Value = "0.03%"
$ws.Range("E51").Value = "0.03%"

# Restore the default cell style so no extra formatting is left behind
$dRange.Style = "Normal"
$eRange.Style = "Normal"

